# Update weekly Ramas de apio price records for rows 2-22
# (Date, Variedad, Calidad, Volumen, Precio min/max/promedio, Precio $/Kg)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44539
$ws.Range("H2").Value = "Americana (o)"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 160
$ws.Range("K2").Value = 6500
$ws.Range("L2").Value = 7000
$ws.Range("M2").Value = 6750
$ws.Range("P2").Value = 6750
$ws.Range("D3").Value = 44414
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 6000
$ws.Range("L3").Value = 7000
$ws.Range("M3").Value = 6500
$ws.Range("P3").Value = 6500
$ws.Range("D4").Value = 45118
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 4000
$ws.Range("L4").Value = 5000
$ws.Range("M4").Value = 4500
$ws.Range("P4").Value = 4500
$ws.Range("D5").Value = 44789
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 80
$ws.Range("K5").Value = 5000
$ws.Range("L5").Value = 6000
$ws.Range("M5").Value = 5500
$ws.Range("P5").Value = 5500
$ws.Range("D6").Value = 44371
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 80
$ws.Range("K6").Value = 7000
$ws.Range("L6").Value = 8000
$ws.Range("M6").Value = 7375
$ws.Range("P6").Value = 7375
$ws.Range("D7").Value = 44764
$ws.Range("H7").Value = "Americana (o)"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = 7500
$ws.Range("P7").Value = 7500
$ws.Range("D8").Value = 44804
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 60
$ws.Range("K8").Value = 5500
$ws.Range("L8").Value = 6000
$ws.Range("M8").Value = 5750
$ws.Range("P8").Value = 5750
$ws.Range("D9").Value = 44253
$ws.Range("H9").Value = "Americana (o)"
$ws.Range("I9").Value = "Segunda"
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 4000
$ws.Range("L9").Value = 4500
$ws.Range("M9").Value = 4250
$ws.Range("P9").Value = 4250
$ws.Range("D10").Value = 44945
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 45
$ws.Range("K10").Value = 6000
$ws.Range("L10").Value = 7000
$ws.Range("M10").Value = 6444
$ws.Range("P10").Value = 6444
$ws.Range("D11").Value = 44410
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 5500
$ws.Range("L11").Value = 6000
$ws.Range("M11").Value = 5750
$ws.Range("P11").Value = 5750
$ws.Range("D12").Value = 44263
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 7000
$ws.Range("L12").Value = 8000
$ws.Range("M12").Value = 7500
$ws.Range("P12").Value = 7500
$ws.Range("D13").Value = 44309
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 50
$ws.Range("K13").Value = 8000
$ws.Range("L13").Value = 9000
$ws.Range("M13").Value = 8500
$ws.Range("P13").Value = 8500
$ws.Range("D14").Value = 44281
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = 5000
$ws.Range("L14").Value = 6000
$ws.Range("M14").Value = 5500
$ws.Range("P14").Value = 5500
$ws.Range("D15").Value = 44259
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 80
$ws.Range("K15").Value = 4000
$ws.Range("L15").Value = 4500
$ws.Range("M15").Value = 4250
$ws.Range("P15").Value = 4250
$ws.Range("D16").Value = 44699
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 50
$ws.Range("K16").Value = 9000
$ws.Range("L16").Value = 9500
$ws.Range("M16").Value = 9250
$ws.Range("P16").Value = 9250
$ws.Range("D17").Value = 45128
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 3500
$ws.Range("L17").Value = 4000
$ws.Range("M17").Value = 3750
$ws.Range("P17").Value = 3750
$ws.Range("D18").Value = 44575
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 160
$ws.Range("K18").Value = 6500
$ws.Range("L18").Value = 7000
$ws.Range("M18").Value = 6750
$ws.Range("P18").Value = 6750
$ws.Range("D19").Value = 44636
$ws.Range("H19").Value = "Americana (o)"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 60
$ws.Range("K19").Value = 8000
$ws.Range("L19").Value = 9000
$ws.Range("M19").Value = 8500
$ws.Range("P19").Value = 8500
$ws.Range("D20").Value = 44497
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 160
$ws.Range("K20").Value = 5000
$ws.Range("L20").Value = 6000
$ws.Range("M20").Value = 5500
$ws.Range("P20").Value = 5500
$ws.Range("D21").Value = 45154
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 100
$ws.Range("K21").Value = 5000
$ws.Range("L21").Value = 6000
$ws.Range("M21").Value = 5500
$ws.Range("P21").Value = 5500
$ws.Range("D22").Value = 44559
$ws.Range("H22").Value = "Americana (o)"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 5000
$ws.Range("L22").Value = 6000
$ws.Range("M22").Value = 5500
$ws.Range("P22").Value = 5500
